$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 1
$ws.Range("G8").Value = 0
$ws.Range("G9").Value = -1
$ws.Range("G10").Value = 15
